$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-03-18 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-03-19 Wednesday", 2)

# Update the math problems table by cell position (row, col) so that the
# replacement that creates a value equal to another old value does not get
# re-matched by a later find/replace.
$tbl = $d.Tables(1)

$tbl.Cell(1,1).Range.Text = "95÷8="
$tbl.Cell(1,2).Range.Text = "57÷5="
$tbl.Cell(1,3).Range.Text = "76÷3="
$tbl.Cell(1,4).Range.Text = "80÷6="
$tbl.Cell(1,5).Range.Text = "76÷7="

$tbl.Cell(5,1).Range.Text = "99÷3="
$tbl.Cell(5,2).Range.Text = "87÷3="
$tbl.Cell(5,3).Range.Text = "64÷2="
$tbl.Cell(5,4).Range.Text = "23÷5="
$tbl.Cell(5,5).Range.Text = "12÷2="

$tbl.Cell(9,1).Range.Text = "62÷3="
$tbl.Cell(9,2).Range.Text = "31÷4="
$tbl.Cell(9,3).Range.Text = "26÷9="
$tbl.Cell(9,4).Range.Text = "47÷8="
$tbl.Cell(9,5).Range.Text = "93÷4="

$tbl.Cell(13,1).Range.Text = "73÷8="
$tbl.Cell(13,2).Range.Text = "89÷4="
$tbl.Cell(13,3).Range.Text = "52÷8="
$tbl.Cell(13,4).Range.Text = "63÷8="
$tbl.Cell(13,5).Range.Text = "67÷8="

$tbl.Cell(17,1).Range.Text = "92÷4="
$tbl.Cell(17,2).Range.Text = "41÷7="
$tbl.Cell(17,3).Range.Text = "33÷7="
$tbl.Cell(17,4).Range.Text = "60÷7="
$tbl.Cell(17,5).Range.Text = "34÷8="
